# [Fonds de solidarite] Add 2020-07-16 data
# Updates "nombre_aides" (column C) and "montant_total" (column D) figures
# for the rows that received new cumulative counts on 2020-07-16.
# Values are written with a leading apostrophe so that Excel keeps storing
# them as text (matching the original inline-string/text cell type used
# throughout this data extract) instead of converting them to numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 3;  C = "828";  D = "2154538.52" },
    @{ Row = 4;  C = "342";  D = "1175844.92" },
    @{ Row = 6;  C = "20";   D = "114500.00" },
    @{ Row = 16; C = "370";  D = "1018417.16" },
    @{ Row = 17; C = "124";  D = "443500.00" },
    @{ Row = 22; C = "290";  D = "755836.59" },
    @{ Row = 25; C = "16";   D = "41500.00" },
    @{ Row = 32; C = "75";   D = "172320.00" },
    @{ Row = 33; C = "425";  D = "1069926.79" },
    @{ Row = 34; C = "155";  D = "556426.95" },
    @{ Row = 35; C = "55";   D = "233974.00" },
    @{ Row = 37; C = "13";   D = "26000.00" },
    @{ Row = 48; C = "61";   D = "152736.00" },
    @{ Row = 49; C = "427";  D = "1127935.90" },
    @{ Row = 50; C = "186";  D = "587500.00" },
    @{ Row = 51; C = "54";   D = "247877.00" },
    @{ Row = 71; C = "187";  D = "433326.09" },
    @{ Row = 72; C = "754";  D = "2027256.67" },
    @{ Row = 73; C = "278";  D = "948266.79" },
    @{ Row = 74; C = "87";   D = "345000.00" },
    @{ Row = 75; C = "18";   D = "93883.20" },
    @{ Row = 76; C = "18";   D = "37500.00" }
)

foreach ($u in $updates) {
    $ws.Range("C$($u.Row)").Formula = "'" + $u.C
    $ws.Range("D$($u.Row)").Formula = "'" + $u.D
}
